# Populate the newly-added 0DTE gamma-wall data row (row 57) and the
# OHLCV/SD columns (C:H) for row 56 that were left blank in the prior save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56: fill in C56:H56 (SPY OPEN/HIGH/LOW/CLOSE/VOLUME/SDs)
$ws.Range("C56").Value = 591.25
$ws.Range("D56").Value = 594.5
$ws.Range("E56").Value = 589.28
$ws.Range("F56").Value = 594.20000000000005
$ws.Range("G56").Value = 74560456
$ws.Range("H56").Value = 0.62288054645945623

# Row 57: new row of data (dates in A:B, gamma-wall metrics in I:II)
$ws.Range("A57").Value = 45793
$ws.Range("B57").Value = 45796
$ws.Range("I57").Value = 0.1724
$ws.Range("J57").Value = 3.7
$ws.Range("K57").Value = 600
$ws.Range("L57").Value = 190809600
$ws.Range("M57").Value = -36207
$ws.Range("N57").Value = 4718
$ws.Range("O57").Value = 40925
$ws.Range("P57").Value = 0.10032964417165027
$ws.Range("Q57").Value = 0.13433636664191603
$ws.Range("R57").Value = 45796
$ws.Range("S57").Value = 0.13433636664191603
$ws.Range("T57").Value = 45807
$ws.Range("U57").Value = 0.11402819405225818
$ws.Range("V57").Value = 45828
$ws.Range("W57").Value = 0.28831155457516516
$ws.Range("X57").Value = 17.333333333333332
$ws.Range("Y57").Value = 605
$ws.Range("Z57").Value = 102444650
$ws.Range("AA57").Value = 8561
$ws.Range("AB57").Value = 1931
$ws.Range("AC57").Value = 10492
$ws.Range("AD57").Value = 0.053866447399864854
$ws.Range("AE57").Value = 0.073522370037266041
$ws.Range("AF57").Value = 45796
$ws.Range("AG57").Value = 0.073522370037266041
$ws.Range("AH57").Value = 45807
$ws.Range("AI57").Value = 0.20507775904065956
$ws.Range("AJ57").Value = 45828
$ws.Range("AK57").Value = 0.3004861241229988
$ws.Range("AL57").Value = 17.333333333333332
$ws.Range("AM57").Value = 595
$ws.Range("AN57").Value = 98653975
$ws.Range("AO57").Value = -9346
$ws.Range("AP57").Value = 24372
$ws.Range("AQ57").Value = 33718
$ws.Range("AR57").Value = 0.051873271616673808
$ws.Range("AS57").Value = 0.10376587581889447
$ws.Range("AT57").Value = 45797
$ws.Range("AU57").Value = 0.11139734906718743
$ws.Range("AV57").Value = 45807
$ws.Range("AW57").Value = 0.28753064361989444
$ws.Range("AX57").Value = 45828
$ws.Range("AY57").Value = 0.11579248671532712
$ws.Range("AZ57").Value = 17.666666666666668
$ws.Range("BA57").Value = 610
$ws.Range("BB57").Value = 78647300
$ws.Range("BC57").Value = -10834
$ws.Range("BD57").Value = 443
$ws.Range("BE57").Value = 11277
$ws.Range("BF57").Value = 0.041353556760566713
$ws.Range("BG57").Value = 0
$ws.Range("BH57").Value = 45828
$ws.Range("BI57").Value = 0.28420571761496805
$ws.Range("BJ57").Value = 45856
$ws.Range("BK57").Value = 0.12248910898955258
$ws.Range("BL57").Value = 45919
$ws.Range("BM57").Value = 0.13396273884085222
$ws.Range("BN57").Value = 74.666666666666671
$ws.Range("BO57").Value = 615
$ws.Range("BP57").Value = 75502320
$ws.Range("BQ57").Value = -538
$ws.Range("BR57").Value = -16
$ws.Range("BS57").Value = 554
$ws.Range("BT57").Value = 0.039699894029095358
$ws.Range("BU57").Value = 0
$ws.Range("BV57").Value = 45828
$ws.Range("BW57").Value = 0.16351495685448858
$ws.Range("BX57").Value = 45856
$ws.Range("BY57").Value = 0.29442681399405263
$ws.Range("BZ57").Value = 45919
$ws.Range("CA57").Value = 0.1780453508118118
$ws.Range("CB57").Value = 74.666666666666671
$ws.Range("CC57").Value = 588
$ws.Range("CD57").Value = -202090308
$ws.Range("CE57").Value = 0.088395544726204381
$ws.Range("CF57").Value = -11757
$ws.Range("CG57").Value = -6810
$ws.Range("CH57").Value = 18567
$ws.Range("CI57").Value = 0.55174524121048407
$ws.Range("CJ57").Value = 45796
$ws.Range("CK57").Value = 0.55174524121048407
$ws.Range("CL57").Value = 45797
$ws.Range("CM57").Value = 0.17197853796479973
$ws.Range("CN57").Value = 45800
$ws.Range("CO57").Value = 0.11169702021694172
$ws.Range("CP57").Value = 4.666666666666667
$ws.Range("CQ57").Value = 585
$ws.Range("CR57").Value = -196755390
$ws.Range("CS57").Value = 0.086062018752857683
$ws.Range("CT57").Value = -18845
$ws.Range("CU57").Value = -7187
$ws.Range("CV57").Value = 26032
$ws.Range("CW57").Value = 0.21121734432160116
$ws.Range("CX57").Value = 45796
$ws.Range("CY57").Value = 0.21121734432160116
$ws.Range("CZ57").Value = 45800
$ws.Range("DA57").Value = 0.12842306682014365
$ws.Range("DB57").Value = 45807
$ws.Range("DC57").Value = 0.19710970684901669
$ws.Range("DD57").Value = 8
$ws.Range("DE57").Value = 575
$ws.Range("DF57").Value = -119524100
$ws.Range("DG57").Value = 0.052280577094322232
$ws.Range("DH57").Value = -38878
$ws.Range("DI57").Value = -6339
$ws.Range("DJ57").Value = 45217
$ws.Range("DK57").Value = 0.024293439173170665
$ws.Range("DL57").Value = 45800
$ws.Range("DM57").Value = 0.10480731543036756
$ws.Range("DN57").Value = 45807
$ws.Range("DO57").Value = 0.30046387378371076
$ws.Range("DP57").Value = 45828
$ws.Range("DQ57").Value = 0.24118772998664317
$ws.Range("DR57").Value = 18.666666666666668
$ws.Range("DS57").Value = 587
$ws.Range("DT57").Value = -96963008
$ws.Range("DU57").Value = 0.042412216574242212
$ws.Range("DV57").Value = -9089
$ws.Range("DW57").Value = -14723
$ws.Range("DX57").Value = 23812
$ws.Range("DY57").Value = 0.44951334737321419
$ws.Range("DZ57").Value = 45796
$ws.Range("EA57").Value = 0.44951334737321419
$ws.Range("EB57").Value = 45797
$ws.Range("EC57").Value = 0.16899083613182309
$ws.Range("ED57").Value = 45800
$ws.Range("EE57").Value = 0.13287608856508623
$ws.Range("EF57").Value = 4.666666666666667
$ws.Range("EG57").Value = 550
$ws.Range("EH57").Value = -90097150
$ws.Range("EI57").Value = 0.039409048020890465
$ws.Range("EJ57").Value = -23977
$ws.Range("EK57").Value = -52142
$ws.Range("EL57").Value = 76119
$ws.Range("EM57").Value = 0
$ws.Range("EN57").Value = 45828
$ws.Range("EO57").Value = 0.63432756446294802
$ws.Range("EP57").Value = 45856
$ws.Range("EQ57").Value = 0.23599187597181914
$ws.Range("ER57").Value = 45884
$ws.Range("ES57").Value = 0.074757266801391817
$ws.Range("ET57").Value = 63
$ws.Range("EU57").Value = 585
$ws.Range("EV57").Value = 351105300
$ws.Range("EW57").Value = -18845
$ws.Range("EX57").Value = -7187
$ws.Range("EY57").Value = 26032
$ws.Range("EZ57").Value = 0.083835402582952231
$ws.Range("FA57").Value = 77174955
$ws.Range("FB57").Value = 0.040579382662681136
$ws.Range("FC57").Value = 0.020815172486981042
$ws.Range("FD57").Value = 45800
$ws.Range("FE57").Value = 0.19159661317586774
$ws.Range("FF57").Value = 45807
$ws.Range("FG57").Value = 0.15667472692403903
$ws.Range("FH57").Value = 45828
$ws.Range("FI57").Value = 0.14306072481674917
$ws.Range("FJ57").Value = 18.666666666666668
$ws.Range("FK57").Value = -273930345
$ws.Range("FL57").Value = 0.11981881913561188
$ws.Range("FM57").Value = 0.21121734432160116
$ws.Range("FN57").Value = 45796
$ws.Range("FO57").Value = 0.21121734432160116
$ws.Range("FP57").Value = 45800
$ws.Range("FQ57").Value = 0.12842306682014365
$ws.Range("FR57").Value = 45807
$ws.Range("FS57").Value = 0.19710970684901669
$ws.Range("FT57").Value = 8
$ws.Range("FU57").Value = 580
$ws.Range("FV57").Value = 333944860
$ws.Range("FW57").Value = -23106
$ws.Range("FX57").Value = -683
$ws.Range("FY57").Value = 23789
$ws.Range("FZ57").Value = 0.079737907056964458
$ws.Range("GA57").Value = 127891160
$ws.Range("GB57").Value = 0.06724648327704473
$ws.Range("GC57").Value = 0.015129114475152153
$ws.Range("GD57").Value = 45814
$ws.Range("GE57").Value = 0.055178637835484483
$ws.Range("GF57").Value = 45828
$ws.Range("GG57").Value = 0.64813924590253147
$ws.Range("GH57").Value = 45919
$ws.Range("GI57").Value = 0.067260160905570021
$ws.Range("GJ57").Value = 60.666666666666664
$ws.Range("GK57").Value = -206053700
$ws.Range("GL57").Value = 0.090129156784450534
$ws.Range("GM57").Value = 0.065686177923521882
$ws.Range("GN57").Value = 45800
$ws.Range("GO57").Value = 0.18529548365304774
$ws.Range("GP57").Value = 45807
$ws.Range("GQ57").Value = 0.22920355227787709
$ws.Range("GR57").Value = 45828
$ws.Range("GS57").Value = 0.094405584563635597
$ws.Range("GT57").Value = 18.666666666666668
$ws.Range("GU57").Value = 590
$ws.Range("GV57").Value = 290131320
$ws.Range("GW57").Value = -27625
$ws.Range("GX57").Value = 25985
$ws.Range("GY57").Value = 53610
$ws.Range("GZ57").Value = 0.069276299771388647
$ws.Range("HA57").Value = 126051140
$ws.Range("HB57").Value = 0.066278981894154557
$ws.Range("HC57").Value = 0.068664987877142564
$ws.Range("HD57").Value = 45807
$ws.Range("HE57").Value = 0.23188826376342173
$ws.Range("HF57").Value = 45828
$ws.Range("HG57").Value = 0.16831581213783547
$ws.Range("HH57").Value = 45838
$ws.Range("HI57").Value = 0.10070396824653867
$ws.Range("HJ57").Value = 31.333333333333332
$ws.Range("HK57").Value = -164080180
$ws.Range("HL57").Value = 0.071769680760116736
$ws.Range("HM57").Value = 0
$ws.Range("HN57").Value = 45799
$ws.Range("HO57").Value = 0.085824625497119764
$ws.Range("HP57").Value = 45800
$ws.Range("HQ57").Value = 0.4226183198970162
$ws.Range("HR57").Value = 45807
$ws.Range("HS57").Value = 0.10378925717902064
$ws.Range("HT57").Value = 9
$ws.Range("HU57").Value = 592
$ws.Range("HV57").Value = -782303
$ws.Range("HW57").Value = -2392173
$ws.Range("HX57").Value = 1901826739
$ws.Range("HY57").Value = -2286204679.5
$ws.Range("HZ57").Value = -384377940.5
$ws.Range("IA57").Value = 0.83187072271059093
$ws.Range("IB57").Value = 4188031418.5
$ws.Range("IC57").Value = 0.12759761654113772
$ws.Range("ID57").Value = 45796
$ws.Range("IE57").Value = 0.12759761654113772
$ws.Range("IF57").Value = 45807
$ws.Range("IG57").Value = 0.13284263927496132
$ws.Range("IH57").Value = 45828
$ws.Range("II57").Value = 0.1924195884587297

# Restore the active-pane selection to match the saved workbook view
# (frozen pane at K2 stays the same; the bottom-right pane selection moves
# from G62 to F44).
$ws.Range("F44").Select()
